$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Clear the previously hard-coded input value in G14 (was 168)
$ws.Range("G14").ClearContents()

# Set the new hard-coded input value in D11 (was blank)
$ws.Range("D11").Value = 24

# Recalculate so dependent formulas (D10,D12,D13,D14,D15,G10,G11,G12,G13,G15) update
$excel.Calculate()

# Update the selected/active cell to D15
$ws.Range("D15").Select()
